$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21, shifting existing rows 21-102 down to 22-103
$ws.Rows(21).Insert()

# Populate the new row 21 with the weekly price-report entry
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = "2022-02-08"
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112043
$ws.Cells.Item(21, 7).Value = "Pepino ensalada"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 7000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 7520
$ws.Cells.Item(21, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(21, 15).Value = "Región del Maule"
$ws.Cells.Item(21, 16).Value = 125
$ws.Cells.Item(21, 17).Value = 60
$ws.Cells.Item(21, 18).Value = "Hortaliza"
